$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 421, shifting existing rows 421-478 down to 422-479.
$ws.Rows.Item(421).Insert()

# Populate the newly inserted row 421 with the new data record.
$ws.Range("A421").Value = 4
$ws.Range("B421").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C421").Value = "Los Lagos"
$ws.Range("D421").Value = 45142
$ws.Range("E421").Value = 10
$ws.Range("F421").Value = 100112043
$ws.Range("G421").Value = "Pepino ensalada"
$ws.Range("H421").Value = "Sin especificar"
$ws.Range("I421").Value = "Primera"
$ws.Range("J421").Value = 400
$ws.Range("K421").Value = 15000
$ws.Range("L421").Value = 15000
$ws.Range("M421").Value = 15000
$ws.Range("N421").Value = "$/caja 60 unidades"
$ws.Range("O421").Value = "Región de Arica y Parinacota"
$ws.Range("P421").Value = 250
$ws.Range("Q421").Value = 60
$ws.Range("R421").Value = "Hortaliza"
